$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H100").Value = 5016.7095
$ws.Range("I100").Value = 2033.1
$ws.Range("J100").Value = 6437.476
$ws.Range("K100").Value = 2033.1
$ws.Range("L100").Value = 6437.476
$ws.Range("M100").Value = -1492.1
$ws.Range("N100").Value = -7519.476
$ws.Range("H133").Value = 31031.25
$ws.Range("J133").Value = 31031.25
$ws.Range("L133").Value = 31031.25
$ws.Range("N133").Value = -41151.25
$ws.Range("H136").Value = 35046
$ws.Range("J136").Value = 35046
$ws.Range("L136").Value = 35046
$ws.Range("N136").Value = -45246

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 18307.244
$ws.Range("I32").Value = 17687.379
$ws.Range("J32").Value = 24041
$ws.Range("K32").Value = 17687.379
$ws.Range("L32").Value = 24041
$ws.Range("M32").Value = -17400.379
$ws.Range("N32").Value = -24615
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = ""
$ws.Range("H63").Value = 3620.4
$ws.Range("I63").Value = 3686.2856
$ws.Range("J63").Value = 3466.6667
$ws.Range("K63").Value = 3686.2856
$ws.Range("L63").Value = 3466.6667
$ws.Range("M63").Value = -3000.2856
$ws.Range("N63").Value = -4838.6667
$ws.Range("H66").Value = 3620.4
$ws.Range("I66").Value = 3686.2856
$ws.Range("J66").Value = 3466.6667
$ws.Range("K66").Value = 18431.428
$ws.Range("L66").Value = 17333.3335
$ws.Range("M66").Value = -14999.428
$ws.Range("N66").Value = -24197.3335
$ws.Range("H102").Value = 2176.7
$ws.Range("I102").Value = 1307.4445
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 1307.4445
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = 314.5554999999999
$ws.Range("N102").Value = -13244

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H82").Value = 19261.625
$ws.Range("I82").Value = 8508.143
$ws.Range("J82").Value = 27625.445
$ws.Range("K82").Value = 8508.143
$ws.Range("L82").Value = 27625.445
$ws.Range("M82").Value = -8125.143
$ws.Range("N82").Value = -28391.445
$ws.Range("H85").Value = 19261.625
$ws.Range("I85").Value = 8508.143
$ws.Range("J85").Value = 27625.445
$ws.Range("K85").Value = 8508.143
$ws.Range("L85").Value = 27625.445
$ws.Range("M85").Value = -7182.143
$ws.Range("N85").Value = -30277.445

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").Value = ""
$ws.Range("H51").Value = 10000
$ws.Range("J51").Value = 10000
$ws.Range("L51").Value = 10000
$ws.Range("N51").Value = -11472
$ws.Range("H59").Value = 73960
$ws.Range("J59").Value = 73960
$ws.Range("L59").Value = 73960
$ws.Range("N59").Value = -76250
$ws.Range("H60").Value = 5593
$ws.Range("I60").Value = 5593
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 5593
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = ""
$ws.Range("N60").Value = -5082
$ws.Range("H61").Value = 10000
$ws.Range("J61").Value = 10000
$ws.Range("L61").Value = 10000
$ws.Range("N61").Value = -10696
$ws.Range("H74").Value = 13981.429
$ws.Range("J74").Value = 13981.429
$ws.Range("L74").Value = 13981.429
$ws.Range("N74").Value = -15729.429
$ws.Range("H77").Value = 13981.429
$ws.Range("J77").Value = 13981.429
$ws.Range("L77").Value = 41944.287
$ws.Range("N77").Value = -50680.287
$ws.Range("H86").Value = 6266.25
$ws.Range("I86").Value = 2798.1428
$ws.Range("J86").Value = 11121.6
$ws.Range("K86").Value = 2798.1428
$ws.Range("L86").Value = 11121.6
$ws.Range("M86").Value = -1675.1428
$ws.Range("N86").Value = -13367.6
$ws.Range("H89").Value = 6266.25
$ws.Range("I89").Value = 2798.1428
$ws.Range("J89").Value = 11121.6
$ws.Range("K89").Value = 13990.714
$ws.Range("L89").Value = 55608
$ws.Range("M89").Value = -8374.714
$ws.Range("N89").Value = -66840

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H2").Value = 4199.1665
$ws.Range("I2").Value = 6254.375
$ws.Range("J2").Value = 88.75
$ws.Range("K2").Value = 37526.25
$ws.Range("L2").Value = 532.5
$ws.Range("M2").Value = -37413.25
$ws.Range("N2").Value = -758.5
$ws.Range("H5").Value = 645.6429000000001
$ws.Range("I5").Value = 528.25
$ws.Range("J5").Value = 1350
$ws.Range("K5").Value = 1584.75
$ws.Range("L5").Value = 4050
$ws.Range("M5").Value = -1472.75
$ws.Range("N5").Value = -4274
$ws.Range("H68").Value = 399.16666
$ws.Range("J68").Value = 403.63635
$ws.Range("L68").Value = 1210.90905
$ws.Range("N68").Value = -2832.90905
$ws.Range("H71").Value = 399.16666
$ws.Range("J71").Value = 403.63635
$ws.Range("L71").Value = 3632.72715
$ws.Range("N71").Value = -11744.72715
$ws.Range("H76").Value = 3500
$ws.Range("J76").Value = 3500
$ws.Range("L76").Value = 10500
$ws.Range("N76").Value = -11266
$ws.Range("H79").Value = 3500
$ws.Range("J79").Value = 3500
$ws.Range("L79").Value = 10500
$ws.Range("N79").Value = -13152
$ws.Range("H131").Value = 1573.119
$ws.Range("I131").Value = 1353.909
$ws.Range("J131").Value = 1650.9032
$ws.Range("K131").Value = 4061.727
$ws.Range("L131").Value = 4952.7096
$ws.Range("M131").Value = 978.2729999999997
$ws.Range("N131").Value = -15032.7096
$ws.Range("H135").Value = 645.6429000000001
$ws.Range("I135").Value = 528.25
$ws.Range("J135").Value = 1350
$ws.Range("K135").Value = 4754.25
$ws.Range("L135").Value = 12150
$ws.Range("M135").Value = -2219.25
$ws.Range("N135").Value = -17220

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H32").Value = 28193.334
$ws.Range("J32").Value = 28193.334
$ws.Range("L32").Value = 28193.334
$ws.Range("N32").Value = -28785.334
$ws.Range("H80").Value = 6559.706
$ws.Range("I80").Value = 9391.5
$ws.Range("J80").Value = 2514.2856
$ws.Range("K80").Value = 9391.5
$ws.Range("L80").Value = 2514.2856
$ws.Range("M80").Value = -8393.5
$ws.Range("N80").Value = -4510.2856
$ws.Range("H83").Value = 6559.706
$ws.Range("I83").Value = 9391.5
$ws.Range("J83").Value = 2514.2856
$ws.Range("K83").Value = 46957.5
$ws.Range("L83").Value = 12571.428
$ws.Range("M83").Value = -41965.5
$ws.Range("N83").Value = -22555.428
$ws.Range("H107").Value = 22732010
$ws.Range("I107").Value = 7271.7856
$ws.Range("J107").Value = 62500304
$ws.Range("K107").Value = 7271.7856
$ws.Range("L107").Value = 62500304
$ws.Range("M107").Value = -5351.7856
$ws.Range("N107").Value = -62504144
$ws.Range("H122").Value = 215325.58
$ws.Range("I122").Value = 273532.53
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 820597.5900000001
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -818147.5900000001
$ws.Range("N122").Value = -10600

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H55").Value = 239.66667
$ws.Range("J55").Value = 227.71428
$ws.Range("L55").Value = 227.71428
$ws.Range("N55").Value = -573.71428
$ws.Range("H132").Value = 8584.093999999999
$ws.Range("I132").Value = 2530.1
$ws.Range("J132").Value = 18674.084
$ws.Range("K132").Value = 7590.299999999999
$ws.Range("L132").Value = 56022.25199999999
$ws.Range("M132").Value = -5060.299999999999
$ws.Range("N132").Value = -61082.25199999999
$ws.Range("H136").Value = 3821.9375
$ws.Range("I136").Value = 2052.6758
$ws.Range("J136").Value = 9773.091
$ws.Range("K136").Value = 6158.0274
$ws.Range("L136").Value = 29319.273
$ws.Range("M136").Value = -3608.0274
$ws.Range("N136").Value = -34419.273

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H8").Value = 300
$ws.Range("J8").Value = 300
$ws.Range("L8").Value = 300
$ws.Range("N8").Value = -580
$ws.Range("H118").Value = 28000
$ws.Range("J118").Value = 28000
$ws.Range("L118").Value = 28000
$ws.Range("N118").Value = -31314
$ws.Range("H126").Value = 1356.6522
$ws.Range("J126").Value = 1903.5454
$ws.Range("L126").Value = 5710.6362
$ws.Range("N126").Value = -10650.6362
